$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "R"/"S" columns (18/19) of the "2024" sheet hold a log of SMS-style
# messages with their timestamps, newest entry first (row 48 down to row
# 195, with row 196 onward being blank). A brand-new message was received,
# so every existing entry needs to shift down one row to make room for it
# at the top of the log.
#
# Walk the range bottom-up so each read happens before it gets overwritten.
for ($i = 196; $i -ge 49; $i--) {
    $srcRow = $i - 1
    $ws.Cells.Item($i, 18).Value2 = $ws.Cells.Item($srcRow, 18).Value2
    $ws.Cells.Item($i, 19).Value2 = $ws.Cells.Item($srcRow, 19).Value2
}

# Record the newly received message at the top of the log.
$ws.Cells.Item(48, 18).Value2 = "balance your axis"
$ws.Cells.Item(48, 19).Value2 = "2024-09-24 12:44:43"

# The "Broadband" category label (column A) was the last row of the fixed
# category list; it now moves down one row to row 205 as well.
$ws.Cells.Item(205, 1).Value2 = $ws.Cells.Item(204, 1).Value2
$ws.Cells.Item(204, 1).Value2 = ""
